{"js": "// Apply the dated worksheet update: refresh the header date and every\n// division-fact answer cell, in document (reading) order. Some answer\n// strings repeat (e.g. \"92\u00f76=15, 2\" appears twice) with different\n// replacements, so matching is done positionally against the ordered\n// list of (old, new) pairs rather than via global find/replace.\nconst replacements = [\n  [\"2026-01-10 Saturday\", \"2026-01-11 Sunday\"],\n  [\"81\u00f74=20, 1\", \"10\u00f75=2, 0\"],\n  [\"66\u00f72=33, 0\", \"63\u00f77=9, 0\"],\n  [\"82\u00f74=20, 2\", \"98\u00f73=32, 2\"],\n  [\"82\u00f73=27, 1\", \"14\u00f74=3, 2\"],\n  [\"31\u00f74=7, 3\", \"58\u00f74=14, 2\"],\n  [\"33\u00f73=11, 0\", \"25\u00f76=4, 1\"],\n  [\"67\u00f74=16, 3\", \"49\u00f76=8, 1\"],\n  [\"60\u00f72=30, 0\", \"13\u00f74=3, 1\"],\n  [\"12\u00f74=3, 0\", \"51\u00f76=8, 3\"],\n  [\"73\u00f76=12, 1\", \"74\u00f73=24, 2\"],\n  [\"33\u00f78=4, 1\", \"60\u00f75=12, 0\"],\n  [\"22\u00f77=3, 1\", \"26\u00f79=2, 8\"],\n  [\"88\u00f72=44, 0\", \"56\u00f74=14, 0\"],\n  [\"85\u00f78=10, 5\", \"44\u00f73=14, 2\"],\n  [\"24\u00f79=2, 6\", \"26\u00f74=6, 2\"],\n  [\"64\u00f73=21, 1\", \"30\u00f79=3, 3\"],\n  [\"63\u00f77=9, 0\", \"93\u00f79=10, 3\"],\n  [\"37\u00f72=18, 1\", \"69\u00f77=9, 6\"],\n  [\"45\u00f77=6, 3\", \"98\u00f73=32, 2\"],\n  [\"92\u00f76=15, 2\", \"36\u00f74=9, 0\"],\n  [\"29\u00f76=4, 5\", \"29\u00f78=3, 5\"],\n  [\"47\u00f75=9, 2\", \"34\u00f78=4, 2\"],\n  [\"92\u00f76=15, 2\", \"78\u00f77=11, 1\"],\n  [\"96\u00f75=19, 1\", \"80\u00f78=10, 0\"],\n  [\"98\u00f77=14, 0\", \"32\u00f73=10, 2\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet next = 0;\nfor (let i = 0; i < paragraphs.items.length && next < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  const [oldText, newText] = replacements[next];\n  if (text === oldText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n    next++;\n  }\n}\n\nawait context.sync();\n\nif (next !== replacements.length) {\n  throw new Error(\n    \"Only matched \" + next + \" of \" + replacements.length + \" expected paragraphs\"\n  );\n}\n", "ps1": "# Apply the dated worksheet update: refresh the header date and every\n# division-fact answer cell, in document order. Some answer strings repeat\n# (e.g. \"92\u00f76=15, 2\" appears twice) with different replacements, so each\n# pair is located with Find/Replace restricted to wdReplaceOne and the\n# search range is advanced past the replacement before looking for the\n# next pair - that keeps duplicate source strings resolved positionally\n# instead of both being overwritten with the same new value.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-10 Saturday\", \"2026-01-11 Sunday\"),\n    @(\"81\u00f74=20, 1\", \"10\u00f75=2, 0\"),\n    @(\"66\u00f72=33, 0\", \"63\u00f77=9, 0\"),\n    @(\"82\u00f74=20, 2\", \"98\u00f73=32, 2\"),\n    @(\"82\u00f73=27, 1\", \"14\u00f74=3, 2\"),\n    @(\"31\u00f74=7, 3\", \"58\u00f74=14, 2\"),\n    @(\"33\u00f73=11, 0\", \"25\u00f76=4, 1\"),\n    @(\"67\u00f74=16, 3\", \"49\u00f76=8, 1\"),\n    @(\"60\u00f72=30, 0\", \"13\u00f74=3, 1\"),\n    @(\"12\u00f74=3, 0\", \"51\u00f76=8, 3\"),\n    @(\"73\u00f76=12, 1\", \"74\u00f73=24, 2\"),\n    @(\"33\u00f78=4, 1\", \"60\u00f75=12, 0\"),\n    @(\"22\u00f77=3, 1\", \"26\u00f79=2, 8\"),\n    @(\"88\u00f72=44, 0\", \"56\u00f74=14, 0\"),\n    @(\"85\u00f78=10, 5\", \"44\u00f73=14, 2\"),\n    @(\"24\u00f79=2, 6\", \"26\u00f74=6, 2\"),\n    @(\"64\u00f73=21, 1\", \"30\u00f79=3, 3\"),\n    @(\"63\u00f77=9, 0\", \"93\u00f79=10, 3\"),\n    @(\"37\u00f72=18, 1\", \"69\u00f77=9, 6\"),\n    @(\"45\u00f77=6, 3\", \"98\u00f73=32, 2\"),\n    @(\"92\u00f76=15, 2\", \"36\u00f74=9, 0\"),\n    @(\"29\u00f76=4, 5\", \"29\u00f78=3, 5\"),\n    @(\"47\u00f75=9, 2\", \"34\u00f78=4, 2\"),\n    @(\"92\u00f76=15, 2\", \"78\u00f77=11, 1\"),\n    @(\"96\u00f75=19, 1\", \"80\u00f78=10, 0\"),\n    @(\"98\u00f77=14, 0\", \"32\u00f73=10, 2\")\n)\n\n$searchStart = 0\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Range($searchStart, $d.Content.End)\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n    if (-not $found) {\n        throw \"Could not find expected text: $oldText\"\n    }\n\n    $searchStart = $rng.End\n}\n"}
